# "Generate Report for Archive"
#
# The report generator re-emitted the localization-status workbook and, in
# doing so, swapped the rows for the "acfd965a-..." and "bf87c32c-..." files
# (bf87c32c now sorts/appears before acfd965a) on every sheet:
#   - Overview (row4 <-> row5 : File Name / zh-cn / de-de status columns)
#   - zh-cn    (row4 <-> row5 : Source File Name / Status / Latest Target File / Latest Handoff Datetime)
#   - de-de    (row4 <-> row5 : Source File Name / Status / Latest Target File / Latest Handoff Datetime)
#
# bf87c32c's status becomes "In Translation" while acfd965a keeps
# "Ready for handoff".  Hyperlink targets keep pointing at the same files
# (acfd965a / bf87c32c) they always did; only which row displays which
# file (and therefore which target) changes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Non-hyperlink status cells for rows 4 & 5 swap content.
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

# Hyperlinked File Name cells (column A) need both their displayed text
# and their link target to be rebuilt consistently, and the individual
# hyperlink entries that were loaded from the original file cannot be
# edited or removed in place - so rebuild the whole collection in the
# correct final order/content.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/20bd3c1928312e4d1fd793077804dbcfba4ef78d/e2e/0d559da7-8441-464a-a2fd-94acb6dc5246.md", "", "", "0d559da7-8441-464a-a2fd-94acb6dc5246.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/20bd3c1928312e4d1fd793077804dbcfba4ef78d/e2e/a3072634-210a-4096-9bd7-fa7631fd8695.md", "", "", "a3072634-210a-4096-9bd7-fa7631fd8695.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6d6197dc0b14af38a6a4721ad262ad2f0017f02b/e2e/bf87c32c-5941-4d78-9481-92095468d7c4.md", "", "", "bf87c32c-5941-4d78-9481-92095468d7c4.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/055921e1375729072e6e047541ad24ccf3d7b68b/e2e/acfd965a-d0df-4166-a4b7-508cad8cf1c4.md", "", "", "acfd965a-d0df-4166-a4b7-508cad8cf1c4.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/055921e1375729072e6e047541ad24ccf3d7b68b/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B4").Value = "In Translation"
$ws.Range("D4").Value = "2016-03-04 08:07:43"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("D5").Value = "2016-03-04 08:08:24"

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/20bd3c1928312e4d1fd793077804dbcfba4ef78d/e2e/0d559da7-8441-464a-a2fd-94acb6dc5246.md", "", "", "0d559da7-8441-464a-a2fd-94acb6dc5246.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/92ae8e9f10b98207a7b7d4061c2a64c5b9320df5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0d559da7-8441-464a-a2fd-94acb6dc5246.51a92bfbfec82bb90d126378087d1d2956804ee4.zh-cn.xlf", "", "", "0d559da7-8441-464a-a2fd-94acb6dc5246.51a92bfbfec82bb90d126378087d1d2956804ee4.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/20bd3c1928312e4d1fd793077804dbcfba4ef78d/e2e/a3072634-210a-4096-9bd7-fa7631fd8695.md", "", "", "a3072634-210a-4096-9bd7-fa7631fd8695.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/92ae8e9f10b98207a7b7d4061c2a64c5b9320df5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a3072634-210a-4096-9bd7-fa7631fd8695.c43b998cfa53e4fd5894af96f9aed119edbde5bc.zh-cn.xlf", "", "", "a3072634-210a-4096-9bd7-fa7631fd8695.c43b998cfa53e4fd5894af96f9aed119edbde5bc.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6d6197dc0b14af38a6a4721ad262ad2f0017f02b/e2e/bf87c32c-5941-4d78-9481-92095468d7c4.md", "", "", "bf87c32c-5941-4d78-9481-92095468d7c4.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4ae28b6feecf817dde92a45cd095577de3db7957/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/bf87c32c-5941-4d78-9481-92095468d7c4.5fead15a35388c8298fa6546c99eb3dcadb757ec.zh-cn.xlf", "", "", "bf87c32c-5941-4d78-9481-92095468d7c4.5fead15a35388c8298fa6546c99eb3dcadb757ec.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/055921e1375729072e6e047541ad24ccf3d7b68b/e2e/acfd965a-d0df-4166-a4b7-508cad8cf1c4.md", "", "", "acfd965a-d0df-4166-a4b7-508cad8cf1c4.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eeff1708c6980c6e536d9a13766db00bdba3d1b7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/acfd965a-d0df-4166-a4b7-508cad8cf1c4.217dce75458345da52833d02167560c3e2aa5950.zh-cn.xlf", "", "", "acfd965a-d0df-4166-a4b7-508cad8cf1c4.217dce75458345da52833d02167560c3e2aa5950.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/055921e1375729072e6e047541ad24ccf3d7b68b/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B4").Value = "In Translation"
$ws.Range("D4").Value = "2016-03-04 08:07:54"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("D5").Value = "2016-03-04 08:08:33"

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/20bd3c1928312e4d1fd793077804dbcfba4ef78d/e2e/0d559da7-8441-464a-a2fd-94acb6dc5246.md", "", "", "0d559da7-8441-464a-a2fd-94acb6dc5246.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7563f0b178a77596849a1505f6e8cdd290ee12bc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0d559da7-8441-464a-a2fd-94acb6dc5246.51a92bfbfec82bb90d126378087d1d2956804ee4.de-de.xlf", "", "", "0d559da7-8441-464a-a2fd-94acb6dc5246.51a92bfbfec82bb90d126378087d1d2956804ee4.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/20bd3c1928312e4d1fd793077804dbcfba4ef78d/e2e/a3072634-210a-4096-9bd7-fa7631fd8695.md", "", "", "a3072634-210a-4096-9bd7-fa7631fd8695.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7563f0b178a77596849a1505f6e8cdd290ee12bc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a3072634-210a-4096-9bd7-fa7631fd8695.c43b998cfa53e4fd5894af96f9aed119edbde5bc.de-de.xlf", "", "", "a3072634-210a-4096-9bd7-fa7631fd8695.c43b998cfa53e4fd5894af96f9aed119edbde5bc.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6d6197dc0b14af38a6a4721ad262ad2f0017f02b/e2e/bf87c32c-5941-4d78-9481-92095468d7c4.md", "", "", "bf87c32c-5941-4d78-9481-92095468d7c4.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/50c58d178a3ba7c3faa5fc2fac6a6b0f5e6a6c3d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/bf87c32c-5941-4d78-9481-92095468d7c4.5fead15a35388c8298fa6546c99eb3dcadb757ec.de-de.xlf", "", "", "bf87c32c-5941-4d78-9481-92095468d7c4.5fead15a35388c8298fa6546c99eb3dcadb757ec.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/055921e1375729072e6e047541ad24ccf3d7b68b/e2e/acfd965a-d0df-4166-a4b7-508cad8cf1c4.md", "", "", "acfd965a-d0df-4166-a4b7-508cad8cf1c4.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4505e56819726e1738036ace43e17e50e21584ee/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/acfd965a-d0df-4166-a4b7-508cad8cf1c4.217dce75458345da52833d02167560c3e2aa5950.de-de.xlf", "", "", "acfd965a-d0df-4166-a4b7-508cad8cf1c4.217dce75458345da52833d02167560c3e2aa5950.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/055921e1375729072e6e047541ad24ccf3d7b68b/.localization-config", "", "", ".localization-config")
